$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match formatting of the existing header cells (bold font, thin border, centered/top alignment)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data values for columns I and J, rows 2-11
$data = @(
    @(8, 8),
    @(6, 7),
    @(8, 9),
    @(8, 8),
    @(7, 7),
    @(7, 8),
    @(5, 5),
    @(6, 7),
    @(9, 9),
    @(4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
